# Plantilla Objetivos - data refresh
# - bumps the report date (B column) from 45931 to 45962 for all existing data rows
# - tweaks a handful of metric cells on specific employee rows
# - resets employee 1014158 (row 9) metrics to 0
# - fills in employee 1025113 (row 42), previously missing its date, and
#   highlights it (along with a brand-new row 43 for employee 1017255) in yellow
# - updates the sheet's scroll/selection view state

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Bump the date column (B) for every existing data row (2-41) to the new
#    reporting date. Row 42 (currently blank) is handled later together with
#    its other fixes.
# ---------------------------------------------------------------------------
$ws.Range("B2:B41").Value = 45962

# ---------------------------------------------------------------------------
# 2) Row-specific metric corrections
# ---------------------------------------------------------------------------
# Employee 1013436 (row 2): flex_max 8 -> 9
$ws.Range("D2").Value = 9

# Employee 1014158 (row 9): all metrics reset to 0
$ws.Range("C9:M9").Value = 0

# Employee 1015856 (row 10): fidepuntos_pospago 10 -> 12
$ws.Range("K10").Value = 12

# Employee 1026014 (row 12): flex_max 6 -> 7, fidepuntos_pospago 8 -> 9
$ws.Range("D12").Value = 7
$ws.Range("K12").Value = 9

# Employee 1016723 (row 13): flex_max 10 -> 12, fidepuntos_pospago 12 -> 14
$ws.Range("D13").Value = 12
$ws.Range("K13").Value = 14

# Employee 1023030 (row 14): every metric collapses to 1
$ws.Range("C14:N14").Value = 1

# ---------------------------------------------------------------------------
# 3) Row 42 (employee 1025113): give it the new date, drop recargas to 1, and
#    apply the new "highlight" style (bold font + thin border + yellow fill)
#    to column A, matching the updated workbook style table.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A42").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A42").Interior.Color = 65535 # RGB(255,255,0)
$ws.Range("A42").Value = 1025113
$ws.Range("B42").Value = 45962
$ws.Range("N42").Value = 1

# ---------------------------------------------------------------------------
# 4) Add new row 43 for employee 1017255, copying row 42's formatting
#    (including the new highlight style) and filling every metric with 1.
# ---------------------------------------------------------------------------
$ws.Range("A42:N42").Copy()
$ws.Range("A43:N43").PasteSpecial(-4122)
$ws.Range("A43").Value = 1017255
$ws.Range("B43").Value = 45962
$ws.Range("C43:N43").Value = 1

# ---------------------------------------------------------------------------
# 5) Update the sheet view: scroll position + active selection
# ---------------------------------------------------------------------------
$ws.Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N44").Select()
